$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove extra columns (J, K) from the mySprTest2 header rows (20, 21)
$ws.Range("J20").ClearContents()
$ws.Range("K20").ClearContents()
$ws.Range("J21").ClearContents()
$ws.Range("K21").ClearContents()

# Row 22: collapse numeric 1,2,3 in I22:K22 into a single text cell I22 = "1,2,3"
$ws.Range("J22").ClearContents()
$ws.Range("K22").ClearContents()
$ws.Range("I22").Value = "1,2,3"

# Update selection to match target state
$ws.Range("O14").Select()
